$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold, border, centered) from existing header cell H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-18
$dataI = @(8, 8, 9, 5, 5, 8, 9, 6, 10, 8, 7, 8, 8, 9, 8, 3, 4)
$dataJ = @(9, 9, 9, 7, 7, 8, 9, 8, 10, 8, 7, 8, 9, 9, 8, 4, 4)

for ($i = 0; $i -lt $dataI.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $dataI[$i]
    $ws.Cells.Item($row, 10).Value = $dataJ[$i]
}
